$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
# Row 65
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
# Row 69
$ws.Range("H69").Value = 52705
$ws.Range("I69").Value = 64693.75
$ws.Range("K69").Value = 194081.25
$ws.Range("M69").Value = -193207.25
# Row 72
$ws.Range("H72").Value = 52705
$ws.Range("I72").Value = 64693.75
$ws.Range("K72").Value = 582243.75
$ws.Range("M72").Value = -577875.75
# Row 88
$ws.Range("H88").Value = 1335
$ws.Range("I88").Value = 1687.6666
$ws.Range("J88").Value = 1142.6364
$ws.Range("K88").Value = 1687.6666
$ws.Range("L88").Value = 1142.6364
$ws.Range("M88").Value = -1281.6666
$ws.Range("N88").Value = -1954.6364
# Row 91
$ws.Range("H91").Value = 1335
$ws.Range("I91").Value = 1687.6666
$ws.Range("J91").Value = 1142.6364
$ws.Range("K91").Value = 1687.6666
$ws.Range("L91").Value = 1142.6364
$ws.Range("M91").Value = -283.6666
$ws.Range("N91").Value = -3950.6364
# Row 107
$ws.Range("H107").Value = 284.93332
$ws.Range("I107").Value = 180.42857
$ws.Range("J107").Value = 376.375
$ws.Range("K107").Value = 180.42857
$ws.Range("L107").Value = 376.375
$ws.Range("M107").Value = 1739.57143
$ws.Range("N107").Value = -4216.375
# Row 111
$ws.Range("H111").Value = 13017.4
$ws.Range("J111").Value = 7793
$ws.Range("L111").Value = 23379
$ws.Range("N111").Value = -29513
# Row 113
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4500
$ws.Range("N113").Value = -11008
$ws.Range("M113").ClearContents()
# Row 137
$ws.Range("H137").Value = 7500
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 10000
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 30000
$ws.Range("M137").Value = -12450
$ws.Range("N137").Value = -35100
# Row 138
$ws.Range("H138").Value = 3019.6428
$ws.Range("I138").Value = 2283
$ws.Range("J138").Value = 3756.2856
$ws.Range("K138").Value = 6849
$ws.Range("L138").Value = 11268.8568
$ws.Range("M138").Value = -1709
$ws.Range("N138").Value = -21548.8568

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4058891.5
$ws.Range("I32").Value = 4121702.2
$ws.Range("J32").Value = 3525000
$ws.Range("K32").Value = 4121702.2
$ws.Range("L32").Value = 3525000
$ws.Range("M32").Value = -4121415.2
$ws.Range("N32").Value = -3525574
# Row 50
$ws.Range("H50").Value = 4015.9
$ws.Range("I50").Value = 5832.5
$ws.Range("K50").Value = 5832.5
$ws.Range("M50").Value = -5118.5
# Row 74
$ws.Range("H74").Value = 3003.4285
$ws.Range("I74").Value = 2984
$ws.Range("J74").Value = 3029.3333
$ws.Range("K74").Value = 2984
$ws.Range("L74").Value = 3029.3333
$ws.Range("M74").Value = -2110
$ws.Range("N74").Value = -4777.3333
# Row 77
$ws.Range("H77").Value = 3003.4285
$ws.Range("I77").Value = 2984
$ws.Range("J77").Value = 3029.3333
$ws.Range("K77").Value = 14920
$ws.Range("L77").Value = 15146.6665
$ws.Range("M77").Value = -10552
$ws.Range("N77").Value = -23882.6665
# Row 135
$ws.Range("H135").Value = 24999
$ws.Range("J135").Value = 24999
$ws.Range("L135").Value = 24999
$ws.Range("N135").Value = -35139

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 6791.6665
$ws.Range("J15").Value = 425
$ws.Range("L15").Value = 425
$ws.Range("N15").Value = -765
# Row 74
$ws.Range("H74").Value = 28250
$ws.Range("J74").Value = 28250
$ws.Range("L74").Value = 28250
$ws.Range("N74").Value = -29998
# Row 77
$ws.Range("H77").Value = 28250
$ws.Range("J77").Value = 28250
$ws.Range("L77").Value = 84750
$ws.Range("N77").Value = -93486
# Row 99
$ws.Range("H99").Value = 2386.625
$ws.Range("I99").Value = 1999
$ws.Range("J99").Value = 2774.25
$ws.Range("K99").Value = 1999
$ws.Range("L99").Value = 2774.25
$ws.Range("M99").Value = -501
$ws.Range("N99").Value = -5770.25
# Row 126
$ws.Range("H126").Value = 2386.625
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 2774.25
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 8322.75
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -13262.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 499
$ws.Range("I5").Value = 499
$ws.Range("K5").Value = 1497
$ws.Range("M5").Value = -1385
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 135
$ws.Range("H135").Value = 499
$ws.Range("I135").Value = 499
$ws.Range("K135").Value = 4491
$ws.Range("M135").Value = -1956

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 9627.5
$ws.Range("I7").Value = 10637.5
$ws.Range("K7").Value = 10637.5
$ws.Range("M7").Value = -10525.5
# Row 55
$ws.Range("H55").Value = 1128.6154
$ws.Range("J55").Value = 1332.3334
$ws.Range("L55").Value = 1332.3334
$ws.Range("N55").Value = -1678.3334
# Row 126
$ws.Range("H126").Value = 9627.5
$ws.Range("I126").Value = 10637.5
$ws.Range("K126").Value = 31912.5
$ws.Range("M126").Value = -29442.5
# Row 132
$ws.Range("H132").Value = 3993.75
$ws.Range("I132").Value = 4158.3335
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 12475.0005
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -9945.000499999998
$ws.Range("N132").Value = -15560

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 10678.111
$ws.Range("I4").Value = 680
$ws.Range("J4").Value = 23175.75
$ws.Range("K4").Value = 680
$ws.Range("L4").Value = 23175.75
$ws.Range("M4").Value = -567
$ws.Range("N4").Value = -23401.75
# Row 45
$ws.Range("H45").Value = 18419.6
$ws.Range("J45").Value = 18419.6
$ws.Range("L45").Value = 18419.6
$ws.Range("N45").Value = -19401.6
# Row 59
$ws.Range("H59").Value = 37995
$ws.Range("J59").Value = 37995
$ws.Range("L59").Value = 37995
$ws.Range("N59").Value = -39471
